$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "9.930.282 €"
$ws.Range("C3").Value = "10.730.507 €"
$ws.Range("C4").Value = "14.510.291 €"
$ws.Range("C5").Value = "12.241.853 €"
$ws.Range("C6").Value = "15.017.875 €"
